$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("S-Matrix")

# 1. Rename "Layer 61 switch" -> "Layer 61 string" and document its meaning in the (new) E29 cell.
$ws.Range("C29").Value = "Layer 61 string"
$ws.Range("E29").Value = "'-1 if L61 is not used, marker-procedure-layers when used."

# 2. Update the main description text in A1 (rewording).
$ws.Range("A1").Value = "The S-matrix contains all sample data. It is a 3-dimensional array defined as a global variable."

# 3. Add a new section header in the (previously empty) row 33, for the WF-properties table,
#    mirroring the existing "For [y] = 4/6 and [z] = i" headers.
$ws.Range("B33").Value = "For [y] = 5 and [z] = i"

# 4. Update the view state: scrolled up a bit and selection moved to B34.
$excel.ActiveWindow.ScrollRow = 25
$ws.Range("B34").Select()

# 5. Column B widened (best-fit) to accommodate the header text.
$ws.Columns.Item(2).ColumnWidth = 19.584
